# Update ring item stats in ItemDataTable sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ItemDataTable")

# Row 12: Ring of Protection - defense 100 -> 200
$ws.Range("D12").Value = 200

# Row 13: Ring of Strength - health 100 -> 200
$ws.Range("E13").Value = 200

# Row 14: Ring of Ciritical - critical 100 -> 200
$ws.Range("F14").Value = 200
